# Updates cryptos list price/volume columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "67.870.43"
Set-TextValue "E2" "  +0.01%  "
Set-TextValue "D3" "3.795.81"
Set-TextValue "E3" "  -1.74%  "
Set-TextValue "E4" "  -0.11%  "
Set-TextValue "D5" "598.02"
Set-TextValue "E5" "  -0.10%  "
Set-TextValue "D6" "168.78"
Set-TextValue "E6" "  -0.49%  "
Set-TextValue "D7" "3.799.13"
Set-TextValue "E7" "  -1.55%  "
Set-TextValue "E8" "  -0.08%  "
Set-TextValue "D9" "0.530"
Set-TextValue "E9" "  +0.19%  "
Set-TextValue "D10" "0.166"
Set-TextValue "E10" "  +1.72%  "
Set-TextValue "D11" "6.50"
Set-TextValue "E11" "  +1.93%  "
Set-TextValue "D12" "0.460"
Set-TextValue "E12" "  +1.22%  "
Set-TextValue "D13" "0.0000273"
Set-TextValue "E13" "  +6.90%  "
Set-TextValue "D14" "36.90"
Set-TextValue "E14" "  +0.37%  "
Set-TextValue "D15" "4.432.45"
Set-TextValue "E15" "  -1.62%  "
Set-TextValue "D16" "3.823.16"
Set-TextValue "E16" "  -1.01%  "
Set-TextValue "D17" "19.10"
Set-TextValue "E17" "  +5.69%  "
Set-TextValue "D18" "67.720.49"
Set-TextValue "E18" "  -0.36%  "
Set-TextValue "E19" "  +0.00%  "
Set-TextValue "E20" "  +0.75%  "
Set-TextValue "D21" "10.61"
Set-TextValue "E21" "  -1.73%  "
Set-TextValue "D22" "468.35"
Set-TextValue "E22" "  +0.64%  "
Set-TextValue "E23" "  -0.84%  "
Set-TextValue "D24" "0.0000151"
Set-TextValue "E24" "  -5.24%  "
Set-TextValue "D25" "83.47"
Set-TextValue "E25" "  +0.52%  "
Set-TextValue "D26" "2.29"
Set-TextValue "E26" "  +3.04%  "
Set-TextValue "D27" "12.20"
Set-TextValue "E27" "  +1.62%  "
Set-TextValue "D28" "10.34"
Set-TextValue "E28" "  +4.07%  "
Set-TextValue "E29" "  +0.07%  "
Set-TextValue "E30" "  -0.32%  "
Set-TextValue "D31" "3.944.72"
Set-TextValue "E31" "  -1.59%  "
Set-TextValue "E32" "  -0.33%  "
Set-TextValue "E33" "  -1.81%  "
Set-TextValue "D34" "30.54"
Set-TextValue "E34" "  -1.78%  "
Set-TextValue "D35" "9.22"
Set-TextValue "E35" "  -1.87%  "
Set-TextValue "D36" "3.757.46"
Set-TextValue "E36" "  -1.72%  "
Set-TextValue "D37" "3.79"
Set-TextValue "E37" "  +2.70%  "
Set-TextValue "E38" "  +1.28%  "
Set-TextValue "D39" "5.92"
Set-TextValue "E39" "  +0.69%  "
Set-TextValue "E40" "  -1.15%  "
Set-TextValue "D41" "0.137"
Set-TextValue "E41" "  -1.89%  "
Set-TextValue "E42" "  +0.08%  "
Set-TextValue "D43" "0.319"
Set-TextValue "E43" "  +2.31%  "
Set-TextValue "E44" "  -0.01%  "
Set-TextValue "E45" "  +2.04%  "
Set-TextValue "E46" "  -0.26%  "
Set-TextValue "D47" "409.10"
Set-TextValue "E47" "  -3.25%  "
Set-TextValue "D48" "46.32"
Set-TextValue "E48" "  -1.64%  "
Set-TextValue "D49" "0.000280"
Set-TextValue "E49" "  -6.48%  "
Set-TextValue "D50" "142.07"
Set-TextValue "E50" "  -0.82%  "
Set-TextValue "E51" "  -0.04%  "
